$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace backslashes with forward-slashes in the train/test path cells
$ws.Range("B2").Value = "../assets/test_cases/test_case21.245/train/"
$ws.Range("C2").Value = "../assets/test_cases/test_case21.245/test/"

# Update the computed total_accuracy value
$ws.Range("AI2").Value = 49.3

# Select C3, matching the author's saved cursor position
$ws.Range("C3").Select()
